# Weekly crime-data refresh: new report period (11/6/2023 - 11/12/2023, Vol 30 No 45)
# and updated crime-complaint figures for the 23rd Precinct (rows 15-30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/number and reporting week ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Cells whose type flips between a numeric value and the literal
#     placeholder text "0" / "***.*" (used when a count, or its % change,
#     is not meaningful) need their style copied from a stable donor cell
#     so the General/Text number format travels with the new value. ---
$zeroDonor = $ws.Range("C14")     # s=14, text "0"
$naDonor   = $ws.Range("E14")     # s=14, text "***.*"
$countDonor = $ws.Range("G14")    # s=15, numeric count style

$zeroDonor.Copy($ws.Range("D15"))
$naDonor.Copy($ws.Range("E15"))
$zeroDonor.Copy($ws.Range("F15"))
$zeroDonor.Copy($ws.Range("D20"))
$naDonor.Copy($ws.Range("E20"))
$zeroDonor.Copy($ws.Range("D26"))
$naDonor.Copy($ws.Range("E26"))
$zeroDonor.Copy($ws.Range("C27"))
$countDonor.Copy($ws.Range("F30"))
$ws.Range("F30").Value = 1

# --- Remaining cells: value-only updates (existing style/format retained) ---
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = -52.173913043478
$ws.Range("N15").Value = -81.034482758620
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 171
$ws.Range("J16").Value = 172
$ws.Range("K16").Value = -0.581395348837
$ws.Range("L16").Value = 6.211180124223
$ws.Range("M16").Value = -21.917808219178
$ws.Range("N16").Value = -78.651685393258
$ws.Range("C17").Value = 11
$ws.Range("E17").Value = 22.222222222222
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 47.058823529411
$ws.Range("I17").Value = 376
$ws.Range("J17").Value = 342
$ws.Range("K17").Value = 9.941520467836
$ws.Range("L17").Value = 16.770186335403
$ws.Range("M17").Value = 95.833333333333
$ws.Range("N17").Value = -28.652751423149
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -6.666666666666
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 168
$ws.Range("K18").Value = -32.142857142857
$ws.Range("L18").Value = 25.274725274725
$ws.Range("M18").Value = 32.558139534883
$ws.Range("N18").Value = -72.195121951219
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -23.913043478260
$ws.Range("I19").Value = 377
$ws.Range("J19").Value = 428
$ws.Range("K19").Value = -11.915887850467
$ws.Range("L19").Value = 18.553459119496
$ws.Range("M19").Value = 66.814159292035
$ws.Range("N19").Value = -26.223091976516
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 71
$ws.Range("K20").Value = -4.054054054054
$ws.Range("L20").Value = 33.962264150943
$ws.Range("M20").Value = 108.823529411765
$ws.Range("N20").Value = -80.810810810810
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 120
$ws.Range("H21").Value = -1.639344262295
$ws.Range("I21").Value = 1126
$ws.Range("J21").Value = 1206
$ws.Range("K21").Value = -6.633499170812
$ws.Range("L21").Value = 16.442605997931
$ws.Range("M21").Value = 43.805874840357
$ws.Range("N21").Value = -58.342582315945
$ws.Range("F22").Value = 1
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 14.285714285714
$ws.Range("F23").Value = 43
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = 30.303030303030
$ws.Range("I23").Value = 363
$ws.Range("J23").Value = 362
$ws.Range("K23").Value = 0.276243093922
$ws.Range("L23").Value = -2.419354838709
$ws.Range("M23").Value = 46.963562753036
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -15
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 8.219178082191
$ws.Range("I24").Value = 808
$ws.Range("J24").Value = 781
$ws.Range("K24").Value = 3.457106274007
$ws.Range("L24").Value = 21.503759398496
$ws.Range("M24").Value = 37.883959044368
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 46.341463414634
$ws.Range("I25").Value = 566
$ws.Range("J25").Value = 480
$ws.Range("K25").Value = 17.916666666666
$ws.Range("L25").Value = 20.425531914893
$ws.Range("M25").Value = -4.391891891891
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = -60
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = -15.517241379310
$ws.Range("L27").Value = -25.757575757575
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = -48.275862068965
$ws.Range("N28").Value = -83.516483516483
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("L29").Value = -48.148148148148
$ws.Range("N29").Value = -83.529411764705
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = 50
